# Add Qty column to BOM. Add SMD attribute to ADC board
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F ("Part") for the new "Qty" column.
$ws.Columns.Item(6).Insert()

# Header
$ws.Cells.Item(1, 6).Value = "Qty"

# Quantities per BOM row (count of reference designators in column A)
$qty = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 2
    10 = 1
    11 = 2
    12 = 3
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 3
    20 = 1
    21 = 1
    22 = 2
}

foreach ($r in $qty.Keys) {
    $ws.Cells.Item($r, 6).Value = $qty[$r]
}

# Summary rows
$ws.Range("E24").Value = "total components"
$ws.Range("F24").Formula = "=SUM(F2:F22)"

$ws.Range("E25").Value = "SMD"
$ws.Range("F25").Formula = "=F24-F22-F2"
